$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Successful")

# Clear the bid-data rows (2-4) that were previously populated with sample
# results from a prior automation run, resetting the sheet back to its blank
# template state (e.g. before sending the "run complete" Outlook notification
# and re-arming the process for the next run). Formatting (e.g. the
# wrap-text style already on column E) is left untouched, only the values.
$ws.Range("A2:I4").ClearContents()

# Reset the active cell/selection to F6 (as saved in the workbook view)
$ws.Range("F6").Select()
